# Duplicate the first worksheet ("20180611") into a brand-new sheet placed
# immediately before it, rename the new sheet to "20180803", and restore the
# per-sheet view/selection state that Excel records after this kind of
# "insert a copy for today's snapshot" edit.

$wb = $excel.ActiveWorkbook

# The existing first sheet (tab order), the one that was active/selected.
$source = $wb.Worksheets.Item(1)

# Copy it, inserting the duplicate immediately before $source -> becomes the
# new sheet #1, pushing the original "20180611" to sheet #2.
$source.Copy($source)

$newSheet = $wb.Worksheets.Item(1)
$original = $wb.Worksheets.Item(2)

# Give the freshly duplicated sheet its new name.
$newSheet.Name = "20180803"

# Restore the view state of the original sheet (now in slot 2): it is no
# longer the selected tab, it scrolled to A4, and the selection is C23.
$original.Activate() | Out-Null
$original.Range("C23").Select() | Out-Null

# Restore/establish the view state of the new sheet (slot 1): it is the
# selected tab, scrolled to A4, with selection G20. Doing this last makes it
# the active sheet on save, matching tabSelected="1".
$newSheet.Activate() | Out-Null
$newSheet.Range("G20").Select() | Out-Null
